$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.382440567016602
$ws.Range("B1").Value = 2.681482315063477
$ws.Range("C1").Value = 3.250740528106689
$ws.Range("D1").Value = 3.245505094528198
$ws.Range("E1").Value = 2.168324708938599
